# ---------------------------------------------------------------------------
# Commit: "Sat, Jul 18, 2020  9:06:34 PM"
#
# Two logical changes:
#
#   1. The single table on slide 5 is switched from the (custom) table
#      style {C61AE62E-69A7-499C-98CF-3C407BE1F3A8} to the built-in table
#      style {82F5E636-4FF8-42BB-AB30-F5580C60C199}.
#
#   2. The deck's theme ("Integral" / "Red Violet" colour scheme) is
#      swapped for the stock PowerPoint "Office Theme" / "Office" colour
#      scheme (the two colour values that previously lived in the two
#      theme parts trade places). We reproduce that by rewriting every
#      slot of the live ThemeColorScheme (the part that actually drives
#      what you see on the slides/slideMaster) with the stock Office RGB
#      values.
# ---------------------------------------------------------------------------

$p = $ppt.ActivePresentation

# --- 1. Table style -------------------------------------------------------

for ($si = 1; $si -le $p.Slides.Count; $si++) {
    $slide = $p.Slides.Item($si)
    for ($shi = 1; $shi -le $slide.Shapes.Count; $shi++) {
        $shape = $slide.Shapes.Item($shi)
        if ($shape.HasTable) {
            $shape.Table.ApplyStyle("{82F5E636-4FF8-42BB-AB30-F5580C60C199}")
        }
    }
}

# --- 2. Theme colours: Red Violet/Integral -> Office -----------------------

$theme = $p.SlideMaster.Theme
$tcs = $theme.ThemeColorScheme

# index : scheme slot : target "Office" RGB
$tcs.Item(1).RGB  = 0         # dk1      000000
$tcs.Item(2).RGB  = 16777215  # lt1      FFFFFF
$tcs.Item(3).RGB  = 6968388   # dk2      44546A
$tcs.Item(4).RGB  = 15132391  # lt2      E7E6E6
$tcs.Item(5).RGB  = 13998939  # accent1  5B9BD5
$tcs.Item(6).RGB  = 3243501   # accent2  ED7D31
$tcs.Item(7).RGB  = 10855845  # accent3  A5A5A5
$tcs.Item(8).RGB  = 49407     # accent4  FFC000
$tcs.Item(9).RGB  = 12874308  # accent5  4472C4
$tcs.Item(10).RGB = 4697456   # accent6  70AD47
$tcs.Item(11).RGB = 12673797  # hlink    0563C1
$tcs.Item(12).RGB = 7491477   # folHlink 954F72
